$d = $word.ActiveDocument

# The document ends with a trailing empty paragraph (just a paragraph
# mark, formatted with sz/szCs 18 in its pPr/rPr) sitting right before
# the section properties. Remove that empty paragraph, merging it away
# so the document body ends with the previous (non-empty) paragraph.

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)

if ($lastPara.Range.Text -eq "`r") {
    # Extend the range one character to the left so the delete also
    # consumes the paragraph mark that separates it from the previous
    # paragraph - this removes the whole empty <w:p> block instead of
    # just clearing its (already empty) contents.
    $delRange = $d.Range($lastPara.Range.Start - 1, $lastPara.Range.End)
    $delRange.Delete()
}
